$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the CONCATENATE formula used to build the C# "repository.Add<HotelModel>(...)"
# snippet so it now reads "repository.Add(...)" (generic type argument removed).
# M2 holds its own (non-shared) formula; M3:M53 share one formula (anchored at M3).
$ws.Range("M2").Formula = '=CONCATENATE("repository.Add(new HotelModel { Name = """,A2,""", Address = """,G2,""", City = """,H2,""", Stars = ",B2," } );")'
$ws.Range("M3:M53").Formula = '=CONCATENATE("repository.Add(new HotelModel { Name = """,A3,""", Address = """,G3,""", City = """,H3,""", Stars = ",B3," } );")'

# Update the sheet view: clear the scrolled-down top-left cell and move the
# selection to M2:M53 (active cell M2) instead of I43.
$ws.Range("M2:M53").Select()
